# Auto-applies the diff: sets R399:R417 to numeric 0, and appends
# new historical-data rows 418-430 (column R left blank, matching the
# "backup" column's inlineStr/empty state for not-yet-backed-up rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: fill R399:R417 (existing empty inlineStr cells) with numeric 0
for ($r = 399; $r -le 417; $r++) {
    $ws.Cells.Item($r, 18).Value = 0
}

# Step 2: append new rows 418-430 of stock data (column R left blank)
# Row 418
$ws.Cells.Item(418, 1).Value = 45517
$ws.Cells.Item(418, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(418, 2).Value = 1479.900024414062
$ws.Cells.Item(418, 3).Value = 1488.949951171875
$ws.Cells.Item(418, 4).Value = 1450
$ws.Cells.Item(418, 5).Value = 1453.300048828125
$ws.Cells.Item(418, 6).Value = 1453.300048828125
$ws.Cells.Item(418, 7).Value = 476761
$ws.Cells.Item(418, 8).Value = 2024
$ws.Cells.Item(418, 9).Value = 8
$ws.Cells.Item(418, 10).Value = 13
$ws.Cells.Item(418, 11).Value = 0
$ws.Cells.Item(418, 12).Value = 0
$ws.Cells.Item(418, 13).Value = 0
$ws.Cells.Item(418, 14).Value = 33
$ws.Cells.Item(418, 15).Value = 0
$ws.Cells.Item(418, 16).Value = 0
$ws.Cells.Item(418, 17).Value = 2

# Row 419
$ws.Cells.Item(419, 1).Value = 45518
$ws.Cells.Item(419, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(419, 2).Value = 1465
$ws.Cells.Item(419, 3).Value = 1472
$ws.Cells.Item(419, 4).Value = 1423
$ws.Cells.Item(419, 5).Value = 1432.900024414062
$ws.Cells.Item(419, 6).Value = 1432.900024414062
$ws.Cells.Item(419, 7).Value = 301988
$ws.Cells.Item(419, 8).Value = 2024
$ws.Cells.Item(419, 9).Value = 8
$ws.Cells.Item(419, 10).Value = 14
$ws.Cells.Item(419, 11).Value = 0
$ws.Cells.Item(419, 12).Value = 0
$ws.Cells.Item(419, 13).Value = 0
$ws.Cells.Item(419, 14).Value = 33
$ws.Cells.Item(419, 15).Value = 0
$ws.Cells.Item(419, 16).Value = 0
$ws.Cells.Item(419, 17).Value = 0

# Row 420
$ws.Cells.Item(420, 1).Value = 45520
$ws.Cells.Item(420, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(420, 2).Value = 1460
$ws.Cells.Item(420, 3).Value = 1503.199951171875
$ws.Cells.Item(420, 4).Value = 1457
$ws.Cells.Item(420, 5).Value = 1500.050048828125
$ws.Cells.Item(420, 6).Value = 1500.050048828125
$ws.Cells.Item(420, 7).Value = 1047796
$ws.Cells.Item(420, 8).Value = 2024
$ws.Cells.Item(420, 9).Value = 8
$ws.Cells.Item(420, 10).Value = 16
$ws.Cells.Item(420, 11).Value = 0
$ws.Cells.Item(420, 12).Value = 0
$ws.Cells.Item(420, 13).Value = 0
$ws.Cells.Item(420, 14).Value = 33
$ws.Cells.Item(420, 15).Value = 0
$ws.Cells.Item(420, 16).Value = 0
$ws.Cells.Item(420, 17).Value = 0

# Row 421
$ws.Cells.Item(421, 1).Value = 45523
$ws.Cells.Item(421, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(421, 2).Value = 1530
$ws.Cells.Item(421, 3).Value = 1562.400024414062
$ws.Cells.Item(421, 4).Value = 1512.900024414062
$ws.Cells.Item(421, 5).Value = 1520.800048828125
$ws.Cells.Item(421, 6).Value = 1520.800048828125
$ws.Cells.Item(421, 7).Value = 1311965
$ws.Cells.Item(421, 8).Value = 2024
$ws.Cells.Item(421, 9).Value = 8
$ws.Cells.Item(421, 10).Value = 19
$ws.Cells.Item(421, 11).Value = 0
$ws.Cells.Item(421, 12).Value = 0
$ws.Cells.Item(421, 13).Value = 0
$ws.Cells.Item(421, 14).Value = 34
$ws.Cells.Item(421, 15).Value = 0
$ws.Cells.Item(421, 16).Value = 0
$ws.Cells.Item(421, 17).Value = 0

# Row 422
$ws.Cells.Item(422, 1).Value = 45524
$ws.Cells.Item(422, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(422, 2).Value = 1530
$ws.Cells.Item(422, 3).Value = 1530.349975585938
$ws.Cells.Item(422, 4).Value = 1487.449951171875
$ws.Cells.Item(422, 5).Value = 1505
$ws.Cells.Item(422, 6).Value = 1505
$ws.Cells.Item(422, 7).Value = 696087
$ws.Cells.Item(422, 8).Value = 2024
$ws.Cells.Item(422, 9).Value = 8
$ws.Cells.Item(422, 10).Value = 20
$ws.Cells.Item(422, 11).Value = 0
$ws.Cells.Item(422, 12).Value = 0
$ws.Cells.Item(422, 13).Value = 0
$ws.Cells.Item(422, 14).Value = 34
$ws.Cells.Item(422, 15).Value = 0
$ws.Cells.Item(422, 16).Value = 0
$ws.Cells.Item(422, 17).Value = 0

# Row 423
$ws.Cells.Item(423, 1).Value = 45525
$ws.Cells.Item(423, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(423, 2).Value = 1500
$ws.Cells.Item(423, 3).Value = 1526.800048828125
$ws.Cells.Item(423, 4).Value = 1497
$ws.Cells.Item(423, 5).Value = 1516.349975585938
$ws.Cells.Item(423, 6).Value = 1516.349975585938
$ws.Cells.Item(423, 7).Value = 259833
$ws.Cells.Item(423, 8).Value = 2024
$ws.Cells.Item(423, 9).Value = 8
$ws.Cells.Item(423, 10).Value = 21
$ws.Cells.Item(423, 11).Value = 0
$ws.Cells.Item(423, 12).Value = 0
$ws.Cells.Item(423, 13).Value = 0
$ws.Cells.Item(423, 14).Value = 34
$ws.Cells.Item(423, 15).Value = 0
$ws.Cells.Item(423, 16).Value = 0
$ws.Cells.Item(423, 17).Value = 0

# Row 424
$ws.Cells.Item(424, 1).Value = 45526
$ws.Cells.Item(424, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(424, 2).Value = 1525
$ws.Cells.Item(424, 3).Value = 1534.800048828125
$ws.Cells.Item(424, 4).Value = 1510.099975585938
$ws.Cells.Item(424, 5).Value = 1515.150024414062
$ws.Cells.Item(424, 6).Value = 1515.150024414062
$ws.Cells.Item(424, 7).Value = 286835
$ws.Cells.Item(424, 8).Value = 2024
$ws.Cells.Item(424, 9).Value = 8
$ws.Cells.Item(424, 10).Value = 22
$ws.Cells.Item(424, 11).Value = 0
$ws.Cells.Item(424, 12).Value = 0
$ws.Cells.Item(424, 13).Value = 0
$ws.Cells.Item(424, 14).Value = 34
$ws.Cells.Item(424, 15).Value = 0
$ws.Cells.Item(424, 16).Value = 0
$ws.Cells.Item(424, 17).Value = 0

# Row 425
$ws.Cells.Item(425, 1).Value = 45527
$ws.Cells.Item(425, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(425, 2).Value = 1519
$ws.Cells.Item(425, 3).Value = 1520
$ws.Cells.Item(425, 4).Value = 1480
$ws.Cells.Item(425, 5).Value = 1485.800048828125
$ws.Cells.Item(425, 6).Value = 1485.800048828125
$ws.Cells.Item(425, 7).Value = 340053
$ws.Cells.Item(425, 8).Value = 2024
$ws.Cells.Item(425, 9).Value = 8
$ws.Cells.Item(425, 10).Value = 23
$ws.Cells.Item(425, 11).Value = 0
$ws.Cells.Item(425, 12).Value = 0
$ws.Cells.Item(425, 13).Value = 0
$ws.Cells.Item(425, 14).Value = 34
$ws.Cells.Item(425, 15).Value = 0
$ws.Cells.Item(425, 16).Value = 0
$ws.Cells.Item(425, 17).Value = 0

# Row 426
$ws.Cells.Item(426, 1).Value = 45530
$ws.Cells.Item(426, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(426, 2).Value = 1498
$ws.Cells.Item(426, 3).Value = 1516.900024414062
$ws.Cells.Item(426, 4).Value = 1489
$ws.Cells.Item(426, 5).Value = 1511.900024414062
$ws.Cells.Item(426, 6).Value = 1511.900024414062
$ws.Cells.Item(426, 7).Value = 278133
$ws.Cells.Item(426, 8).Value = 2024
$ws.Cells.Item(426, 9).Value = 8
$ws.Cells.Item(426, 10).Value = 26
$ws.Cells.Item(426, 11).Value = 0
$ws.Cells.Item(426, 12).Value = 0
$ws.Cells.Item(426, 13).Value = 0
$ws.Cells.Item(426, 14).Value = 35
$ws.Cells.Item(426, 15).Value = 0
$ws.Cells.Item(426, 16).Value = 0
$ws.Cells.Item(426, 17).Value = 0

# Row 427
$ws.Cells.Item(427, 1).Value = 45532
$ws.Cells.Item(427, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(427, 2).Value = 1522
$ws.Cells.Item(427, 3).Value = 1552
$ws.Cells.Item(427, 4).Value = 1513.550048828125
$ws.Cells.Item(427, 5).Value = 1519.099975585938
$ws.Cells.Item(427, 6).Value = 1519.099975585938
$ws.Cells.Item(427, 7).Value = 644917
$ws.Cells.Item(427, 8).Value = 2024
$ws.Cells.Item(427, 9).Value = 8
$ws.Cells.Item(427, 10).Value = 28
$ws.Cells.Item(427, 11).Value = 0
$ws.Cells.Item(427, 12).Value = 0
$ws.Cells.Item(427, 13).Value = 0
$ws.Cells.Item(427, 14).Value = 35
$ws.Cells.Item(427, 15).Value = 0
$ws.Cells.Item(427, 16).Value = 0
$ws.Cells.Item(427, 17).Value = 0

# Row 428
$ws.Cells.Item(428, 1).Value = 45533
$ws.Cells.Item(428, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(428, 2).Value = 1519.099975585938
$ws.Cells.Item(428, 3).Value = 1524.599975585938
$ws.Cells.Item(428, 4).Value = 1496.099975585938
$ws.Cells.Item(428, 5).Value = 1509.300048828125
$ws.Cells.Item(428, 6).Value = 1509.300048828125
$ws.Cells.Item(428, 7).Value = 240045
$ws.Cells.Item(428, 8).Value = 2024
$ws.Cells.Item(428, 9).Value = 8
$ws.Cells.Item(428, 10).Value = 29
$ws.Cells.Item(428, 11).Value = 0
$ws.Cells.Item(428, 12).Value = 0
$ws.Cells.Item(428, 13).Value = 0
$ws.Cells.Item(428, 14).Value = 35
$ws.Cells.Item(428, 15).Value = 0
$ws.Cells.Item(428, 16).Value = 0
$ws.Cells.Item(428, 17).Value = 0

# Row 429
$ws.Cells.Item(429, 1).Value = 45534
$ws.Cells.Item(429, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(429, 2).Value = 1516.199951171875
$ws.Cells.Item(429, 3).Value = 1530
$ws.Cells.Item(429, 4).Value = 1510
$ws.Cells.Item(429, 5).Value = 1514.550048828125
$ws.Cells.Item(429, 6).Value = 1514.550048828125
$ws.Cells.Item(429, 7).Value = 172436
$ws.Cells.Item(429, 8).Value = 2024
$ws.Cells.Item(429, 9).Value = 8
$ws.Cells.Item(429, 10).Value = 30
$ws.Cells.Item(429, 11).Value = 0
$ws.Cells.Item(429, 12).Value = 0
$ws.Cells.Item(429, 13).Value = 0
$ws.Cells.Item(429, 14).Value = 35
$ws.Cells.Item(429, 15).Value = 0
$ws.Cells.Item(429, 16).Value = 0
$ws.Cells.Item(429, 17).Value = 0

# Row 430
$ws.Cells.Item(430, 1).Value = 45537
$ws.Cells.Item(430, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(430, 2).Value = 1524.900024414062
$ws.Cells.Item(430, 3).Value = 1524.900024414062
$ws.Cells.Item(430, 4).Value = 1502
$ws.Cells.Item(430, 5).Value = 1513.5
$ws.Cells.Item(430, 6).Value = 1513.5
$ws.Cells.Item(430, 7).Value = 189679
$ws.Cells.Item(430, 8).Value = 2024
$ws.Cells.Item(430, 9).Value = 9
$ws.Cells.Item(430, 10).Value = 2
$ws.Cells.Item(430, 11).Value = 0
$ws.Cells.Item(430, 12).Value = 0
$ws.Cells.Item(430, 13).Value = 0
$ws.Cells.Item(430, 14).Value = 36
$ws.Cells.Item(430, 15).Value = 0
$ws.Cells.Item(430, 16).Value = 0
$ws.Cells.Item(430, 17).Value = 0

